$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 19 de Septiembre de 2020 a las 01:37"

# Row 4
$ws.Range("B4").Value = 6920923
$ws.Range("C4").Value = 46327
$ws.Range("D4").Value = 4186008
$ws.Range("E4").Value = 2531812
$ws.Range("G4").Value = 890
$ws.Range("H4").Value = 203103

# Row 6
$ws.Range("B6").Value = 4497434
$ws.Range("C6").Value = 39991
$ws.Range("D6").Value = 3789139
$ws.Range("E6").Value = 572438
$ws.Range("G6").Value = 826
$ws.Range("H6").Value = 135857

# Row 29
$ws.Range("B29").Value = 141911
$ws.Range("C29").Value = 1044
$ws.Range("D29").Value = 123723
$ws.Range("E29").Value = 8983
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = 9205

# Row 31
$ws.Range("A31").Value = "Ecuador"
$ws.Range("B31").Value = 124129
$ws.Range("C31").Value = 1872
$ws.Range("D31").Value = 97063
$ws.Range("E31").Value = 16022
$ws.Range("G31").Value = 15
$ws.Range("H31").Value = 11044

# Row 32
$ws.Range("A32").Value = "Catar"
$ws.Range("B32").Value = 122917
$ws.Range("C32").Value = 224
$ws.Range("D32").Value = 119822
$ws.Range("E32").Value = 2886
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = 209

# Row 37
$ws.Range("B37").Value = 101772
$ws.Range("C37").Value = 131
$ws.Range("D37").Value = 87958
$ws.Range("E37").Value = 8081
$ws.Range("G37").Value = 18
$ws.Range("H37").Value = 5733

# Row 47
$ws.Range("A47").Value = "Japon"
$ws.Range("B47").Value = 77494
$ws.Range("C47").Value = 485
$ws.Range("D47").Value = 69899
$ws.Range("E47").Value = 6113
$ws.Range("G47").Value = 9
$ws.Range("H47").Value = 1482

# Row 48
$ws.Range("A48").Value = "Polonia"
$ws.Range("B48").Value = 77328
$ws.Range("C48").Value = 757
$ws.Range("D48").Value = 63312
$ws.Range("E48").Value = 11746
$ws.Range("G48").Value = 17
$ws.Range("H48").Value = 2270

# Row 54
$ws.Range("B54").Value = 63879
$ws.Range("C54").Value = 690
$ws.Range("D54").Value = 56700
$ws.Range("E54").Value = 6959

# Row 58
$ws.Range("B58").Value = 56956
$ws.Range("C58").Value = 221
$ws.Range("D58").Value = 48305
$ws.Range("E58").Value = 7557
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 1094

# Row 63
$ws.Range("A63").Value = "Chequia"
$ws.Range("B63").Value = 46262
$ws.Range("C63").Value = 2107
$ws.Range("D63").Value = 23858
$ws.Range("E63").Value = 21909
$ws.Range("G63").Value = 6
$ws.Range("H63").Value = 495

# Row 64
$ws.Range("A64").Value = "Ghana"
$ws.Range("B64").Value = 45760
$ws.Range("C64").Value = 46
$ws.Range("D64").Value = 44973
$ws.Range("E64").Value = 492
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 295

# Row 65
$ws.Range("A65").Value = "Moldavia"
$ws.Range("B65").Value = 45648
$ws.Range("C65").Value = 665
$ws.Range("D65").Value = 33734
$ws.Range("E65").Value = 10728
$ws.Range("G65").Value = 16
$ws.Range("H65").Value = 1186

# Row 66
$ws.Range("A66").Value = "Kirguistan"
$ws.Range("B66").Value = 45244
$ws.Range("C66").Value = 91
$ws.Range("D66").Value = 41415
$ws.Range("E66").Value = 2766
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 1063

# Row 93
$ws.Range("B93").Value = 12769
$ws.Range("C93").Value = 125
$ws.Range("E93").Value = 2131

# Row 96
$ws.Range("A96").Value = "Guinea"
$ws.Range("B96").Value = 10231
$ws.Range("C96").Value = 77
$ws.Range("D96").Value = 9660
$ws.Range("E96").Value = 508
$ws.Range("H96").Value = 63

# Row 97
$ws.Range("A97").Value = "Namibia"
$ws.Range("B97").Value = 10207
$ws.Range("C97").Value = 129
$ws.Range("D97").Value = 7947
$ws.Range("E97").Value = 2152
$ws.Range("H97").Value = 108

# Row 100
$ws.Range("B100").Value = 9568
$ws.Range("C100").Value = 74
$ws.Range("D100").Value = 8077
$ws.Range("E100").Value = 1458

# Row 103
$ws.Range("B103").Value = 8696
$ws.Range("C103").Value = 18
$ws.Range("D103").Value = 7848
$ws.Range("E103").Value = 795

# Row 107
$ws.Range("A107").Value = "Montenegro"
$ws.Range("B107").Value = 7711
$ws.Range("C107").Value = 208
$ws.Range("D107").Value = 4997
$ws.Range("E107").Value = 2581
$ws.Range("G107").Value = 4
$ws.Range("H107").Value = 133

# Row 108
$ws.Range("A108").Value = "Zimbabue"
$ws.Range("B108").Value = 7647
$ws.Range("C108").Value = 14
$ws.Range("D108").Value = 5883
$ws.Range("E108").Value = 1540
$ws.Range("H108").Value = 224

# Row 109
$ws.Range("B109").Value = 7361
$ws.Range("C109").Value = 15
$ws.Range("D109").Value = 6882
$ws.Range("E109").Value = 318

# Row 120
$ws.Range("A120").Value = "Congo"
$ws.Range("B120").Value = 4980
$ws.Range("C120").Value = 46
$ws.Range("D120").Value = 3887
$ws.Range("E120").Value = 1004
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = 89

# Row 121
$ws.Range("A121").Value = "Nicaragua"
$ws.Range("B121").Value = 4961
$ws.Range("D121").Value = 2913
$ws.Range("E121").Value = 1901
$ws.Range("H121").Value = 147

# Row 122
$ws.Range("B122").Value = 4786
$ws.Range("C122").Value = 4
$ws.Range("E122").Value = 2894

# Row 123
$ws.Range("A123").Value = "Surinam"
$ws.Range("B123").Value = 4691
$ws.Range("C123").Value = 20
$ws.Range("D123").Value = 4280
$ws.Range("E123").Value = 315
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 96

# Row 124
$ws.Range("A124").Value = "Ruanda"
$ws.Range("C124").Value = 18
$ws.Range("D124").Value = 2845
$ws.Range("E124").Value = 1801
$ws.Range("G124").Value = 2
$ws.Range("H124").Value = 25

# Row 131
$ws.Range("A131").Value = "Trinidad yTobago"
$ws.Range("B131").Value = 3651
$ws.Range("C131").Value = 217
$ws.Range("D131").Value = 1586
$ws.Range("E131").Value = 2005
$ws.Range("H131").Value = 60

# Row 132
$ws.Range("A132").Value = "Lituania"
$ws.Range("B132").Value = 3565
$ws.Range("C132").Value = 61
$ws.Range("D132").Value = 2181
$ws.Range("E132").Value = 1297
$ws.Range("H132").Value = 87

# Row 133
$ws.Range("A133").Value = "Mayotte"
$ws.Range("B133").Value = 3541
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 2964
$ws.Range("E133").Value = 537
$ws.Range("H133").Value = 40

# Row 135
$ws.Range("B135").Value = 3485
$ws.Range("C135").Value = 12
$ws.Range("D135").Value = 1973
$ws.Range("E135").Value = 1404
$ws.Range("G135").Value = 1
$ws.Range("H135").Value = 108

# Row 148
$ws.Range("A148").Value = "Guinea-Bisau"
$ws.Range("B148").Value = 2303
$ws.Range("C148").Value = 28
$ws.Range("D148").Value = 1127
$ws.Range("E148").Value = 1137
$ws.Range("H148").Value = 39

# Row 149
$ws.Range("A149").Value = "Benin"
$ws.Range("B149").Value = 2280
$ws.Range("D149").Value = 1950
$ws.Range("E149").Value = 290
$ws.Range("H149").Value = 40

# Row 154
$ws.Range("B154").Value = 1890
$ws.Range("C154").Value = 14
$ws.Range("D154").Value = 1603
$ws.Range("E154").Value = 241

# Row 214
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# Row 215
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
